$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Update existing property values in place -----------------------------
$ws.Range("B3").Value2 = "0.1.7"
$ws.Range("B6").Value2 = "draft"
$ws.Range("B8").Value2 = "2024-11-22T12:33:30-06:00"
$ws.Range("B10").Value2 = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws.Range("B11").Value2 = "Bob Milius (bmilius@nmdp.org)"

# --- Make room for a new "Jurisdiction" row at row 12 ----------------------
# Read the four rows that need to shift down (Description, Purpose,
# Copyright, Immutable) BEFORE writing anything, then write them back out
# one row lower, working from the bottom up so we never overwrite a cell
# we still need to read.
$a12 = $ws.Range("A12").Value2
$b12 = $ws.Range("B12").Value2
$a13 = $ws.Range("A13").Value2
$b13 = $ws.Range("B13").Value2
$a14 = $ws.Range("A14").Value2
$b14 = $ws.Range("B14").Value2
$a15 = $ws.Range("A15").Value2
$b15 = $ws.Range("B15").Value2

$ws.Range("A16").Value2 = $a15
$ws.Range("B16").Value2 = $b15
$ws.Range("A16:B16").Style = $ws.Range("A15:B15").Style

$ws.Range("A15").Value2 = $a14
$ws.Range("B15").Value2 = $b14

$ws.Range("A14").Value2 = $a13
$ws.Range("B14").Value2 = $b13

$ws.Range("A13").Value2 = $a12
$ws.Range("B13").Value2 = $b12

# New Jurisdiction row
$ws.Range("A12").Value2 = "Jurisdiction"
$ws.Range("B12").Value2 = ""
